$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column C header date, matching style of B1 (bold, bordered, centered)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "13-01-2023"

# Rewrite data rows 2-53: funds sorted alphabetically (rows 2-51),
# then avg (row 52) and total (row 53) at the bottom; new column C values added.
$ws.Cells.Item(2, 1).Value = "1810 Renta variable"
$ws.Cells.Item(2, 2).Value = 579804.47
$ws.Cells.Item(2, 3).Value = 615125.18
$ws.Cells.Item(3, 1).Value = "1822 Raices Valores Negociables"
$ws.Cells.Item(3, 2).Value = 1125849.18
$ws.Cells.Item(3, 3).Value = 1125857.31
$ws.Cells.Item(4, 1).Value = "Adcap IOL Acciones Argentina"
$ws.Cells.Item(4, 2).Value = 117739.51
$ws.Cells.Item(4, 3).Value = 127982.88
$ws.Cells.Item(5, 1).Value = "Allaria Acciones"
$ws.Cells.Item(5, 2).Value = 212982.03
$ws.Cells.Item(5, 3).Value = 212925.38
$ws.Cells.Item(6, 1).Value = "Alpha Acciones"
$ws.Cells.Item(6, 2).Value = 505817.74
$ws.Cells.Item(6, 3).Value = 505965.04
$ws.Cells.Item(7, 1).Value = "Alpha Mega"
$ws.Cells.Item(7, 2).Value = 1247657.65
$ws.Cells.Item(7, 3).Value = 1248017.56
$ws.Cells.Item(8, 1).Value = "Alpha Recursos Naturales"
$ws.Cells.Item(8, 2).Value = 487183.81
$ws.Cells.Item(8, 3).Value = 542205.6
$ws.Cells.Item(9, 1).Value = "Alpha planeam equil"
$ws.Cells.Item(9, 2).Value = 26848.79
$ws.Cells.Item(9, 3).Value = 26952.41
$ws.Cells.Item(10, 1).Value = "Alpha renta balan global"
$ws.Cells.Item(10, 2).Value = 898923.75
$ws.Cells.Item(10, 3).Value = 901326.49
$ws.Cells.Item(11, 1).Value = "Argenfunds"
$ws.Cells.Item(11, 2).Value = 37353.64
$ws.Cells.Item(11, 3).Value = 37364.61
$ws.Cells.Item(12, 1).Value = "Arpenta acciones"
$ws.Cells.Item(12, 2).Value = 11805.46
$ws.Cells.Item(12, 3).Value = 11801.98
$ws.Cells.Item(13, 1).Value = "Arpenta ex Mercosur"
$ws.Cells.Item(13, 2).Value = 28991.39
$ws.Cells.Item(13, 3).Value = 29005.93
$ws.Cells.Item(14, 1).Value = "Balanz"
$ws.Cells.Item(14, 2).Value = 1092459.64
$ws.Cells.Item(14, 3).Value = 1092300.83
$ws.Cells.Item(15, 1).Value = "Bull Market"
$ws.Cells.Item(15, 2).Value = 234651.57
$ws.Cells.Item(15, 3).Value = 255256.87
$ws.Cells.Item(16, 1).Value = "CMA acciones"
$ws.Cells.Item(16, 2).Value = 540524.31
$ws.Cells.Item(16, 3).Value = 507065.21
$ws.Cells.Item(17, 1).Value = "Compass Crecimiento"
$ws.Cells.Item(17, 2).Value = 2111027.88
$ws.Cells.Item(17, 3).Value = 2090732.48
$ws.Cells.Item(18, 1).Value = "Consultatio Acciones Argentina"
$ws.Cells.Item(18, 2).Value = 2224279.86
$ws.Cells.Item(18, 3).Value = 2224566.25
$ws.Cells.Item(19, 1).Value = "Consultatio Renta Variable"
$ws.Cells.Item(19, 2).Value = 995074.59
$ws.Cells.Item(19, 3).Value = 995746.24
$ws.Cells.Item(20, 1).Value = "Delta Acciones"
$ws.Cells.Item(20, 2).Value = 237529.61
$ws.Cells.Item(20, 3).Value = 237215.71
$ws.Cells.Item(21, 1).Value = "Delta Internacional"
$ws.Cells.Item(21, 2).Value = 7507.43
$ws.Cells.Item(21, 3).Value = 7498.84
$ws.Cells.Item(22, 1).Value = "Delta Latinoamerica"
$ws.Cells.Item(22, 2).Value = 16001.67
$ws.Cells.Item(22, 3).Value = 16011.47
$ws.Cells.Item(23, 1).Value = "Delta Recursos Naturales"
$ws.Cells.Item(23, 2).Value = 1849616.67
$ws.Cells.Item(23, 3).Value = 1848703.78
$ws.Cells.Item(24, 1).Value = "Delta Select"
$ws.Cells.Item(24, 2).Value = 2056160.31
$ws.Cells.Item(24, 3).Value = 2057737.32
$ws.Cells.Item(25, 1).Value = "Delta gestion V"
$ws.Cells.Item(25, 2).Value = 697652.69
$ws.Cells.Item(25, 3).Value = 697857.3
$ws.Cells.Item(26, 1).Value = "FBA Acciones Argentinas"
$ws.Cells.Item(26, 2).Value = 837093.1
$ws.Cells.Item(26, 3).Value = 857810.35
$ws.Cells.Item(27, 1).Value = "FBA Calificado"
$ws.Cells.Item(27, 2).Value = 823407.99
$ws.Cells.Item(27, 3).Value = 841844.76
$ws.Cells.Item(28, 1).Value = "Fima Acciones"
$ws.Cells.Item(28, 2).Value = 1343412.19
$ws.Cells.Item(28, 3).Value = 1459212.93
$ws.Cells.Item(29, 1).Value = "Fima PB Acciones"
$ws.Cells.Item(29, 2).Value = 622840.87
$ws.Cells.Item(29, 3).Value = 652445.66
$ws.Cells.Item(30, 1).Value = "Gainvest Renta Variable"
$ws.Cells.Item(30, 2).Value = 99230.79
$ws.Cells.Item(30, 3).Value = 98788.89
$ws.Cells.Item(31, 1).Value = "Galileo Acciones"
$ws.Cells.Item(31, 2).Value = 6080452.74
$ws.Cells.Item(31, 3).Value = 6197094.41
$ws.Cells.Item(32, 1).Value = "Goal Acciones Argentinas"
$ws.Cells.Item(32, 2).Value = 132530.34
$ws.Cells.Item(32, 3).Value = 132550.47
$ws.Cells.Item(33, 1).Value = "Goal acciones plus"
$ws.Cells.Item(33, 2).Value = 24222.87
$ws.Cells.Item(33, 3).Value = 24177.25
$ws.Cells.Item(34, 1).Value = "HF Acciones Argentinas"
$ws.Cells.Item(34, 2).Value = 501205.92
$ws.Cells.Item(34, 3).Value = 491351.58
$ws.Cells.Item(35, 1).Value = "HF Acciones Lideres"
$ws.Cells.Item(35, 2).Value = 890426.16
$ws.Cells.Item(35, 3).Value = 891680.06
$ws.Cells.Item(36, 1).Value = "IAM Renta Variable"
$ws.Cells.Item(36, 2).Value = 120674.03
$ws.Cells.Item(36, 3).Value = 126296.71
$ws.Cells.Item(37, 1).Value = "IEB Value"
$ws.Cells.Item(37, 2).Value = 35433.55
$ws.Cells.Item(37, 3).Value = 35438.87
$ws.Cells.Item(38, 1).Value = "Lombardi"
$ws.Cells.Item(38, 2).Value = 146018.28
$ws.Cells.Item(38, 3).Value = 168968.18
$ws.Cells.Item(39, 1).Value = "MAF"
$ws.Cells.Item(39, 2).Value = 101939.35
$ws.Cells.Item(39, 3).Value = 102030.77
$ws.Cells.Item(40, 1).Value = "Megainver"
$ws.Cells.Item(40, 2).Value = 109322.8
$ws.Cells.Item(40, 3).Value = 109287.81
$ws.Cells.Item(41, 1).Value = "Pellegrini Acciones"
$ws.Cells.Item(41, 2).Value = 323614.49
$ws.Cells.Item(41, 3).Value = 323961.47
$ws.Cells.Item(42, 1).Value = "Pionero Acciones"
$ws.Cells.Item(42, 2).Value = 752391.42
$ws.Cells.Item(42, 3).Value = 752047.52
$ws.Cells.Item(43, 1).Value = "Premier Renta Variable"
$ws.Cells.Item(43, 2).Value = 182433.71
$ws.Cells.Item(43, 3).Value = 222467.76
$ws.Cells.Item(44, 1).Value = "Quinquela Acciones"
$ws.Cells.Item(44, 2).Value = 352786.47
$ws.Cells.Item(44, 3).Value = 352905.89
$ws.Cells.Item(45, 1).Value = "Rofex 20 Renta Variable"
$ws.Cells.Item(45, 2).Value = 245597.84
$ws.Cells.Item(45, 3).Value = 245914.59
$ws.Cells.Item(46, 1).Value = "SBS Acciones Argentina"
$ws.Cells.Item(46, 2).Value = 1281351.12
$ws.Cells.Item(46, 3).Value = 1281210.23
$ws.Cells.Item(47, 1).Value = "Schroeder RV"
$ws.Cells.Item(47, 2).Value = 2237061.58
$ws.Cells.Item(47, 3).Value = 2238299.03
$ws.Cells.Item(48, 1).Value = "Supefondo RV"
$ws.Cells.Item(48, 2).Value = 3815299.8
$ws.Cells.Item(48, 3).Value = 3981705.87
$ws.Cells.Item(49, 1).Value = "Superfondo "
$ws.Cells.Item(49, 2).Value = 3751419.83
$ws.Cells.Item(49, 3).Value = 3753808.9
$ws.Cells.Item(50, 1).Value = "Supergestion"
$ws.Cells.Item(50, 2).Value = 844626.27
$ws.Cells.Item(50, 3).Value = 844640.45
$ws.Cells.Item(51, 1).Value = "Toronto Trust Multimercado"
$ws.Cells.Item(51, 2).Value = 158599.93
$ws.Cells.Item(51, 3).Value = 149476.03
$ws.Cells.Item(52, 1).Value = "avg"
$ws.Cells.Item(52, 2).Value = 863136.74
$ws.Cells.Item(52, 3).Value = 875012.78
$ws.Cells.Item(53, 1).Value = "total"
$ws.Cells.Item(53, 2).Value = 43156837.09
$ws.Cells.Item(53, 3).Value = 43750639.11

Write-Host "done"
